# Generate Report for Handback
# Updates the handoff/handback timestamp cells across the Overview, zh-cn,
# and de-de sheets to reflect a newly generated report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview!G2 - "Latest HO Xliff Generate Date"
$overview.Range("G2").Value = "2016-10-25 02:54:47"

# zh-cn!H2 - "Correspond Handoff Datetime"
$zhcn.Range("H2").Value = "2016-10-25 02:54:34"
# zh-cn!K2 - "Correspond Handback DateTime"
$zhcn.Range("K2").Value = "2016-10-25 02:55:12"

# de-de!H2 - "Correspond Handoff Datetime" (shares value with Overview!G2)
$dede.Range("H2").Value = "2016-10-25 02:54:47"
# de-de!K2 - "Correspond Handback DateTime"
$dede.Range("K2").Value = "2016-10-25 02:55:29"
